$d = $word.ActiveDocument

$d.Content.Find.Execute("405÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "990÷5=", 2)
$d.Content.Find.Execute("852÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "517÷8=", 2)
$d.Content.Find.Execute("146÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "922÷7=", 2)
$d.Content.Find.Execute("796÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "394÷8=", 2)
$d.Content.Find.Execute("321÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "383÷8=", 2)
$d.Content.Find.Execute("358÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "354÷9=", 2)
$d.Content.Find.Execute("407÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "287÷9=", 2)
$d.Content.Find.Execute("760÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "749÷2=", 2)
$d.Content.Find.Execute("364÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "636÷7=", 2)
$d.Content.Find.Execute("564÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "319÷9=", 2)
$d.Content.Find.Execute("724÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "181÷2=", 2)
$d.Content.Find.Execute("463÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "965÷5=", 2)
$d.Content.Find.Execute("351÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "648÷5=", 2)
$d.Content.Find.Execute("467÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "938÷7=", 2)
$d.Content.Find.Execute("250÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "428÷6=", 2)
$d.Content.Find.Execute("637÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "157÷7=", 2)
$d.Content.Find.Execute("560÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "340÷7=", 2)
$d.Content.Find.Execute("372÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "193÷5=", 2)
$d.Content.Find.Execute("209÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "314÷6=", 2)
$d.Content.Find.Execute("389÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "452÷2=", 2)
$d.Content.Find.Execute("797÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "537÷6=", 2)
$d.Content.Find.Execute("371÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "777÷9=", 2)
$d.Content.Find.Execute("178÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "550÷2=", 2)
$d.Content.Find.Execute("730÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "354÷7=", 2)
$d.Content.Find.Execute("114÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "265÷3=", 2)
